$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 495, shifting existing rows 495:534 down to 496:535
$ws.Rows("495:495").Insert()

# Populate the new row 495 with the new weekly price observation
$ws.Range("A495").Value = 9
$ws.Range("B495").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C495").Value = "Metropolitana"
$ws.Range("D495").Value = 45013
$ws.Range("E495").Value = 13
$ws.Range("F495").Value = 100112044
$ws.Range("G495").Value = "Perejil"
$ws.Range("H495").Value = "Sin especificar"
$ws.Range("I495").Value = "Primera"
$ws.Range("J495").Value = 70
$ws.Range("K495").Value = 12000
$ws.Range("L495").Value = 13000
$ws.Range("M495").Value = 12500
$ws.Range("N495").Value = "$/docena de atados"
$ws.Range("O495").Value = "Región Metropolitana"
$ws.Range("P495").Value = 4167
$ws.Range("Q495").Value = 3
$ws.Range("R495").Value = "Hortaliza"
